$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without Excel auto-converting a
# date-looking string ("YYYY-MM-DD") into a real date serial number. We briefly
# flip the cell to Text format while we type the value, then restore the original
# "General" look by pasting-in just the number format from a donor cell that
# already carries the normal General/right-aligned style used throughout the sheet.
function Set-LiteralText($cellAddr, $text) {
    $donor = $ws.Range("D13")
    $target = $ws.Range($cellAddr)
    $target.NumberFormat = "@"
    $target.Value = $text
    $donor.Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# --- Remove the "23.90.0.2" bad-driver row (old row 4). ---
# Deleting the entire row shifts every row below it up by one, which is exactly what
# the target layout needs: the Totals row becomes row 4, and the whole "Good Drivers"
# table (old rows 11-20) becomes rows 10-19 (sheet now ends at row 24 instead of 25).
$ws.Rows.Item(4).Delete()

# --- Bad Drivers table updates ---
$ws.Range("C3").Value = 175
$ws.Range("D3").Value = 94.59999999999999

# Totals row (now row 4) reflects the single remaining bad driver
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 175

# --- Good Drivers table updates (rows now 10-19 after the shift) ---

# Row 12: Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3
$ws.Range("B12").Value = 34181
$ws.Range("D12").Value = 99.90000000000001

# Row 13: Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5
Set-LiteralText "E13" "2024-08-13"

# Row 14: Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1 -> 22.130.0.5
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5"
$ws.Range("B14").Value = 18738
$ws.Range("D14").Value = 99.90000000000001
Set-LiteralText "E14" "2024-01-20"

# Row 15: Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1 -> 23.20.1.1
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1"
$ws.Range("B15").Value = 13533
Set-LiteralText "E15" "2023-12-19"

# Row 16: Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3 -> 22.170.2.1
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1"
$ws.Range("B16").Value = 19083
Set-LiteralText "E16" "2022-11-22"

# Row 17: Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5 -> 22.100.0.3
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3"
$ws.Range("B17").Value = 12988
$ws.Range("D17").Value = 100
Set-LiteralText "E17" "2022-05-01"

# Row 18: Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1 (vintage date only changes)
Set-LiteralText "E18" "2022-05-01"

# Row 19: Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6 stays as-is (no further edits needed)

# --- Nudge the sheet's recognized extent out to column J / row 24 ---
# The template reserves columns A-J (see <cols>) and a block through row 24, even
# though F:J and rows 20-24 hold no values. Touching a cell property (without
# altering its value) on the far corner makes Excel register it as part of the
# sheet's used range/dimension, matching the original report's A1:J24 extent,
# while keeping the cell itself empty and unstyled.
$corner = $ws.Range("J24")
$corner.FormulaHidden = $corner.FormulaHidden()
$corner.Style = "Normal"
